$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 227, shifting existing rows 227-328 down to 228-329.
$ws.Rows("227:227").Insert()

# Populate the newly inserted row 227 with the new weekly record.
$ws.Cells.Item(227, 1).Value = 4
$ws.Cells.Item(227, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(227, 3).Value = "Los Lagos"
$ws.Cells.Item(227, 4).Value = 44917
$ws.Cells.Item(227, 5).Value = 10
$ws.Cells.Item(227, 6).Value = 100112044
$ws.Cells.Item(227, 7).Value = "Perejil"
$ws.Cells.Item(227, 8).Value = "Sin especificar"
$ws.Cells.Item(227, 9).Value = "Primera"
$ws.Cells.Item(227, 10).Value = 80
$ws.Cells.Item(227, 11).Value = 5000
$ws.Cells.Item(227, 12).Value = 5000
$ws.Cells.Item(227, 13).Value = 5000
$ws.Cells.Item(227, 14).Value = "$/docena de atados (2 kilos)"
$ws.Cells.Item(227, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(227, 16).Value = 2500
$ws.Cells.Item(227, 17).Value = 2
$ws.Cells.Item(227, 18).Value = "Hortaliza"
